$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A2 from numeric 400000002 to the text value "UK1312"
$ws.Range("A2").Value = "UK1312"
